# end of day dump, working on finding articles for pathways
#
# Change all the ">0.05" FDR placeholder values in column E of the
# PubMed_ClusterONE sheet (both the " >0.05" with a stray leading space and
# the plain ">0.05" variant) to ">0.1", and move the active cell selection
# to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PubMed_ClusterONE")

$rows = @(6,7,8,11,12,13,14,15,16,18,19,20,23,24,25,26,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,44,45,46,47,48,49,50,51,52,53)

foreach ($r in $rows) {
    $ws.Range("E$r").Value = ">0.1"
}

$ws.Range("C2").Select()
